# Scheduled-runner style market data refresh across the per-job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Only currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns (H-N) on specific rows are refreshed
# with newer market-board values; everything else is left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6999.923
$ws.Range("I18").Value = 4181.727
$ws.Range("K18").Value = 4181.727
$ws.Range("M18").Value = -3897.727
$ws.Range("H82").Value = 222
$ws.Range("I82").Value = 222
$ws.Range("K82").Value = 666
$ws.Range("M82").Value = -260
$ws.Range("H85").Value = 222
$ws.Range("I85").Value = 222
$ws.Range("K85").Value = 666
$ws.Range("M85").Value = 738
$ws.Range("H104").Value = 90
$ws.Range("I104").Value = 90
$ws.Range("K104").Value = 270
$ws.Range("M104").Value = 1477
$ws.Range("H115").Value = 970.5714
$ws.Range("J115").Value = 998.75
$ws.Range("L115").Value = 2996.25
$ws.Range("N115").Value = -6130.25
$ws.Range("H127").Value = 1972.6
$ws.Range("I127").Value = 1217.75
$ws.Range("K127").Value = 3653.25
$ws.Range("M127").Value = 1306.75
$ws.Range("H135").Value = 2114.1052
$ws.Range("I135").Value = 2114.1052
$ws.Range("K135").Value = 19026.9468
$ws.Range("M135").Value = -16491.9468

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20637.566
$ws.Range("I32").Value = 21995.143
$ws.Range("K32").Value = 21995.143
$ws.Range("M32").Value = -21708.143
$ws.Range("H61").Value = 6411.1665
$ws.Range("I61").Value = 1809.3684
$ws.Range("J61").Value = 23898
$ws.Range("K61").Value = 1809.3684
$ws.Range("L61").Value = 23898
$ws.Range("M61").Value = -1597.3684
$ws.Range("N61").Value = -24322
$ws.Range("H74").Value = 235968.23
$ws.Range("I74").Value = 462073.06
$ws.Range("J74").Value = 9863.385
$ws.Range("K74").Value = 462073.06
$ws.Range("L74").Value = 9863.385
$ws.Range("M74").Value = -461199.06
$ws.Range("N74").Value = -11611.385
$ws.Range("H77").Value = 235968.23
$ws.Range("I77").Value = 462073.06
$ws.Range("J77").Value = 9863.385
$ws.Range("K77").Value = 2310365.3
$ws.Range("L77").Value = 49316.925
$ws.Range("M77").Value = -2305997.3
$ws.Range("N77").Value = -58052.925
$ws.Range("H102").Value = 1867.55
$ws.Range("I102").Value = 1824.4445
$ws.Range("K102").Value = 1824.4445
$ws.Range("M102").Value = -202.4445000000001
$ws.Range("H110").Value = 26664.4
$ws.Range("I110").Value = 29425.773
$ws.Range("K110").Value = 29425.773
$ws.Range("M110").Value = -27380.773
$ws.Range("H131").Value = 103112.4
$ws.Range("J131").Value = 103112.4
$ws.Range("L131").Value = 103112.4
$ws.Range("N131").Value = -113192.4
$ws.Range("H132").Value = 1152.7646
$ws.Range("I132").Value = 974.8182
$ws.Range("J132").Value = 2271.2856
$ws.Range("K132").Value = 2924.4546
$ws.Range("L132").Value = 6813.8568
$ws.Range("M132").Value = -394.4546
$ws.Range("N132").Value = -11873.8568
$ws.Range("H136").Value = 6411.1665
$ws.Range("I136").Value = 1809.3684
$ws.Range("J136").Value = 23898
$ws.Range("K136").Value = 5428.1052
$ws.Range("L136").Value = 71694
$ws.Range("M136").Value = -2878.1052
$ws.Range("N136").Value = -76794

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 119995
$ws.Range("J122").Value = 119995
$ws.Range("L122").Value = 119995
$ws.Range("N122").Value = -129795
$ws.Range("H134").Value = 1675.92
$ws.Range("I134").Value = 1122.381
$ws.Range("K134").Value = 3367.143
$ws.Range("M134").Value = -832.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1735.2759
$ws.Range("I58").Value = 1445.2858
$ws.Range("J58").Value = 2496.5
$ws.Range("K58").Value = 1445.2858
$ws.Range("L58").Value = 2496.5
$ws.Range("M58").Value = -1242.2858
$ws.Range("N58").Value = -2902.5
$ws.Range("H132").Value = 143994.42
$ws.Range("I132").Value = 167827
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 503481
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -500951
$ws.Range("N132").Value = -8057
$ws.Range("H136").Value = 1735.2759
$ws.Range("I136").Value = 1445.2858
$ws.Range("J136").Value = 2496.5
$ws.Range("K136").Value = 4335.857400000001
$ws.Range("L136").Value = 7489.5
$ws.Range("M136").Value = -1785.857400000001
$ws.Range("N136").Value = -12589.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 504.44446
$ws.Range("I33").Value = 220
$ws.Range("K33").Value = 1320
$ws.Range("M33").Value = -1037
$ws.Range("H56").Value = 6414.8335
$ws.Range("I56").Value = 6414.8335
$ws.Range("K56").Value = 6414.8335
$ws.Range("M56").Value = -5884.8335
$ws.Range("H87").Value = 8286.571
$ws.Range("I87").Value = 7601.2
$ws.Range("K87").Value = 22803.6
$ws.Range("M87").Value = -21555.6
$ws.Range("H90").Value = 8286.571
$ws.Range("I90").Value = 7601.2
$ws.Range("K90").Value = 68410.8
$ws.Range("M90").Value = -62170.8
$ws.Range("H141").Value = 6269.143
$ws.Range("I141").Value = 5376.8
$ws.Range("K141").Value = 16130.4
$ws.Range("M141").Value = -10950.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7274.875
$ws.Range("I80").Value = 2199.6667
$ws.Range("J80").Value = 10320
$ws.Range("K80").Value = 2199.6667
$ws.Range("L80").Value = 10320
$ws.Range("M80").Value = -1201.6667
$ws.Range("N80").Value = -12316
$ws.Range("H83").Value = 7274.875
$ws.Range("I83").Value = 2199.6667
$ws.Range("J83").Value = 10320
$ws.Range("K83").Value = 10998.3335
$ws.Range("L83").Value = 51600
$ws.Range("M83").Value = -6006.333500000001
$ws.Range("N83").Value = -61584
$ws.Range("H97").Value = 1778.4166
$ws.Range("I97").Value = 1951.6666
$ws.Range("J97").Value = 1605.1666
$ws.Range("K97").Value = 1951.6666
$ws.Range("L97").Value = 1605.1666
$ws.Range("M97").Value = -1455.6666
$ws.Range("N97").Value = -2597.1666
$ws.Range("H102").Value = 2295.9
$ws.Range("I102").Value = 1378.6923
$ws.Range("K102").Value = 1378.6923
$ws.Range("M102").Value = 243.3077000000001
$ws.Range("H132").Value = 2529.7673
$ws.Range("I132").Value = 2366.7144
$ws.Range("K132").Value = 7100.1432
$ws.Range("M132").Value = -4570.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4111.1665
$ws.Range("I7").Value = 4694.125
$ws.Range("K7").Value = 4694.125
$ws.Range("M7").Value = -4582.125
$ws.Range("H68").Value = 2536.6667
$ws.Range("I68").Value = 2536.6667
$ws.Range("K68").Value = 2536.6667
$ws.Range("M68").Value = -1787.6667
$ws.Range("H71").Value = 2536.6667
$ws.Range("I71").Value = 2536.6667
$ws.Range("K71").Value = 12683.3335
$ws.Range("M71").Value = -8939.333500000001
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H122").Value = 3535.25
$ws.Range("I122").Value = 3450.9375
$ws.Range("K122").Value = 10352.8125
$ws.Range("M122").Value = -7902.8125
$ws.Range("H126").Value = 4111.1665
$ws.Range("I126").Value = 4694.125
$ws.Range("K126").Value = 14082.375
$ws.Range("M126").Value = -11612.375
$ws.Range("H132").Value = 4040.647
$ws.Range("I132").Value = 3692.2856
$ws.Range("J132").Value = 5666.3335
$ws.Range("K132").Value = 11076.8568
$ws.Range("L132").Value = 16999.0005
$ws.Range("M132").Value = -8546.856800000001
$ws.Range("N132").Value = -22059.0005
$ws.Range("H136").Value = 4028.5862
$ws.Range("I136").Value = 3558.261
$ws.Range("J136").Value = 5831.5
$ws.Range("K136").Value = 10674.783
$ws.Range("L136").Value = 17494.5
$ws.Range("M136").Value = -8124.782999999999
$ws.Range("N136").Value = -22594.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 11197.8
$ws.Range("I3").Value = 7994.5
$ws.Range("K3").Value = 7994.5
$ws.Range("M3").Value = -7880.5
$ws.Range("H81").Value = 9533.6
$ws.Range("I81").Value = 11204.5
$ws.Range("J81").Value = 2850
$ws.Range("K81").Value = 22409
$ws.Range("L81").Value = 5700
$ws.Range("M81").Value = -21348
$ws.Range("N81").Value = -7822
$ws.Range("H84").Value = 9533.6
$ws.Range("I84").Value = 11204.5
$ws.Range("J84").Value = 2850
$ws.Range("K84").Value = 112045
$ws.Range("L84").Value = 28500
$ws.Range("M84").Value = -106741
$ws.Range("N84").Value = -39108
$ws.Range("H107").Value = 1066
$ws.Range("I107").Value = 1254.8572
$ws.Range("K107").Value = 3764.5716
$ws.Range("M107").Value = -1844.5716
$ws.Range("H113").Value = 1626
$ws.Range("I113").Value = 1410.6666
$ws.Range("K113").Value = 4231.9998
$ws.Range("M113").Value = -2061.9998
$ws.Range("H122").Value = 65150.434
$ws.Range("I122").Value = 75059.875
$ws.Range("J122").Value = 1730
$ws.Range("K122").Value = 225179.625
$ws.Range("L122").Value = 5190
$ws.Range("M122").Value = -222729.625
$ws.Range("N122").Value = -10090
$ws.Range("H126").Value = 3331.95
$ws.Range("I126").Value = 2684.4119
$ws.Range("K126").Value = 8053.2357
$ws.Range("M126").Value = -5583.2357
$ws.Range("H132").Value = 49248.6
$ws.Range("I132").Value = 52627.07
$ws.Range("K132").Value = 157881.21
$ws.Range("M132").Value = -155351.21
$ws.Range("H136").Value = 15786.782
$ws.Range("I136").Value = 17639.2
$ws.Range("J136").Value = 3437.3333
$ws.Range("K136").Value = 52917.60000000001
$ws.Range("L136").Value = 10311.9999
$ws.Range("M136").Value = -50367.60000000001
$ws.Range("N136").Value = -15411.9999
$ws.Range("H141").Value = 120939.3
$ws.Range("I141").Value = 120998
$ws.Range("J141").Value = 120935.625
$ws.Range("K141").Value = 120998
$ws.Range("L141").Value = 120935.625
$ws.Range("M141").Value = -115818
$ws.Range("N141").Value = -131295.625
